# Update the "Estado de Cuenta" worker rows:
# - Swap the worker identity between row 16 and row 17
# - Update the "Periodo Mora" value from 2506 to 2507 on both rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: now LUIS ANTONIO DIAZ MORALES / 9286921
$ws.Range("C16").Value = "9286921"
$ws.Range("D16").Value = "LUIS ANTONIO DIAZ MORALES"
$ws.Range("E16").Value = "2507"

# Row 17: now BEATRIZ MARIA SIERRA CAMPO / 30774023
$ws.Range("C17").Value = "30774023"
$ws.Range("D17").Value = "BEATRIZ MARIA SIERRA CAMPO"
$ws.Range("E17").Value = "2507"
